$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 14 (current row 14 is blank, row 15 starts the
# "Breathing" block). This pushes everything from row 14 onward down by 2 rows,
# matching the diff where old row 15 -> new row 17, etc.
$ws.Rows.Item(14).Resize(2).EntireRow.Insert() | Out-Null

# Fill the two newly inserted rows in column E with the new method signatures,
# using the same direct formatting as the other plain text cells in that
# column (E3:E13) -- copy/paste-formats reuses the existing cell style
# instead of minting a brand new one.
$ws.Range("E14").Value = "DisplayStartMessage(): string"
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null

$ws.Range("E15").Value = "DisplayEndMessage(): string"
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column E widens slightly (best-fit) to accommodate the new, longer entries.
$ws.Columns.Item(5).ColumnWidth = 28.3

# Update the active selection to match the post-edit state.
$ws.Range("E16").Select() | Out-Null
